$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.145036666666667
$ws.Range("H2").Value = 3.43511
$ws.Range("I2").Value = 0.4953865629219574
$ws.Range("J2").Value = 0.4953865629219574
$ws.Range("M2").Value = 4.993165333333334
$ws.Range("N2").Value = 14.979496
$ws.Range("O2").Value = 0.06779298131037136
$ws.Range("P2").Value = 0.06779298131037137
$ws.Range("Q2").Value = 5.717357389395556
$ws.Range("R2").Value = 51.45621650456
$ws.Range("S2").Value = 0.03358373200157736
$ws.Range("T2").Value = 0.03358373200157737

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.145036666666667
$ws.Range("H3").Value = 3.43511
$ws.Range("I3").Value = 0.4953865629219574
$ws.Range("J3").Value = 0.4953865629219574
$ws.Range("O3").Value = 0.5355771637189464
$ws.Range("P3").Value = 0.5355771637189464
$ws.Range("Q3").Value = 45.16818696261667
$ws.Range("R3").Value = 406.51368266355
$ws.Range("S3").Value = 0.2653177303142193
$ws.Range("T3").Value = 0.2653177303142193

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.145036666666667
$ws.Range("H4").Value = 3.43511
$ws.Range("I4").Value = 0.4953865629219574
$ws.Range("J4").Value = 0.4953865629219574
$ws.Range("M4").Value = 29.08216166666666
$ws.Range("N4").Value = 87.24648499999999
$ws.Range("O4").Value = 0.3948530262300277
$ws.Range("P4").Value = 0.3948530262300277
$ws.Range("Q4").Value = 33.30014145426111
$ws.Range("R4").Value = 299.70127308835
$ws.Range("S4").Value = 0.1956048835234269
$ws.Range("T4").Value = 0.1956048835234269

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.145036666666667
$ws.Range("H5").Value = 3.43511
$ws.Range("I5").Value = 0.4953865629219574
$ws.Range("J5").Value = 0.4953865629219574
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.130869
$ws.Range("N5").Value = 0.392607
$ws.Range("O5").Value = 0.001776828740654623
$ws.Range("P5").Value = 0.001776828740654624
$ws.Range("Q5").Value = 0.14984980353
$ws.Range("R5").Value = 1.34864823177
$ws.Range("S5").Value = 0.000880217082733844
$ws.Range("T5").Value = 0.0008802170827338441

# Row 6
$ws.Range("G6").Value = 0.4713496666666666
$ws.Range("I6").Value = 0.2039238551060172
$ws.Range("J6").Value = 0.2039238551060172
$ws.Range("M6").Value = 4.993165333333334
$ws.Range("N6").Value = 14.979496
$ws.Range("O6").Value = 0.06779298131037136
$ws.Range("P6").Value = 0.06779298131037137
$ws.Range("Q6").Value = 2.353526815478222
$ws.Range("R6").Value = 21.181741339304
$ws.Range("S6").Value = 0.0138246060979411
$ws.Range("T6").Value = 0.01382460609794111

# Row 7
$ws.Range("G7").Value = 0.4713496666666666
$ws.Range("I7").Value = 0.2039238551060172
$ws.Range("J7").Value = 0.2039238551060172
$ws.Range("O7").Value = 0.5355771637189464
$ws.Range("P7").Value = 0.5355771637189464
$ws.Range("S7").Value = 0.1092169599323141
$ws.Range("T7").Value = 0.1092169599323141

# Row 8
$ws.Range("G8").Value = 0.4713496666666666
$ws.Range("I8").Value = 0.2039238551060172
$ws.Range("J8").Value = 0.2039238551060172
$ws.Range("M8").Value = 29.08216166666666
$ws.Range("N8").Value = 87.24648499999999
$ws.Range("O8").Value = 0.3948530262300277
$ws.Range("P8").Value = 0.3948530262300277
$ws.Range("Q8").Value = 13.70786720752944
$ws.Range("R8").Value = 123.370804867765
$ws.Range("S8").Value = 0.08051995130910458
$ws.Range("T8").Value = 0.0805199513091046

# Row 9
$ws.Range("G9").Value = 0.4713496666666666
$ws.Range("I9").Value = 0.2039238551060172
$ws.Range("J9").Value = 0.2039238551060172
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.130869
$ws.Range("N9").Value = 0.392607
$ws.Range("O9").Value = 0.001776828740654623
$ws.Range("P9").Value = 0.001776828740654624
$ws.Range("Q9").Value = 0.061685059527
$ws.Range("R9").Value = 0.555165535743
$ws.Range("S9").Value = 0.0003623377666574605
$ws.Range("T9").Value = 0.0003623377666574606

# Row 10
$ws.Range("G10").Value = 0.6323219999999999
$ws.Range("H10").Value = 1.896966
$ws.Range("I10").Value = 0.2735666300991275
$ws.Range("J10").Value = 0.2735666300991275
$ws.Range("M10").Value = 4.993165333333334
$ws.Range("N10").Value = 14.979496
$ws.Range("O10").Value = 0.06779298131037136
$ws.Range("P10").Value = 0.06779298131037137
$ws.Range("Q10").Value = 3.157288289904
$ws.Range("R10").Value = 28.415594609136
$ws.Range("S10").Value = 0.01854589744145142
$ws.Range("T10").Value = 0.01854589744145143

# Row 11
$ws.Range("G11").Value = 0.6323219999999999
$ws.Range("H11").Value = 1.896966
$ws.Range("I11").Value = 0.2735666300991275
$ws.Range("J11").Value = 0.2735666300991275
$ws.Range("O11").Value = 0.5355771637189464
$ws.Range("P11").Value = 0.5355771637189464
$ws.Range("Q11").Value = 24.94316483307
$ws.Range("R11").Value = 224.48848349763
$ws.Range("S11").Value = 0.1465160398366408
$ws.Range("T11").Value = 0.1465160398366409

# Row 12
$ws.Range("G12").Value = 0.6323219999999999
$ws.Range("H12").Value = 1.896966
$ws.Range("I12").Value = 0.2735666300991275
$ws.Range("J12").Value = 0.2735666300991275
$ws.Range("M12").Value = 29.08216166666666
$ws.Range("N12").Value = 87.24648499999999
$ws.Range("O12").Value = 0.3948530262300277
$ws.Range("P12").Value = 0.3948530262300277
$ws.Range("Q12").Value = 18.38929062939
$ws.Range("R12").Value = 165.50361566451
$ws.Range("S12").Value = 0.1080186117701911
$ws.Range("T12").Value = 0.1080186117701911

# Row 13
$ws.Range("G13").Value = 0.6323219999999999
$ws.Range("H13").Value = 1.896966
$ws.Range("I13").Value = 0.2735666300991275
$ws.Range("J13").Value = 0.2735666300991275
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.130869
$ws.Range("N13").Value = 0.392607
$ws.Range("O13").Value = 0.001776828740654623
$ws.Range("P13").Value = 0.001776828740654624
$ws.Range("Q13").Value = 0.082751347818
$ws.Range("R13").Value = 0.744762130362
$ws.Range("S13").Value = 0.0004860810508441619
$ws.Range("T13").Value = 0.000486081050844162

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.062692
$ws.Range("H14").Value = 0.188076
$ws.Range("I14").Value = 0.02712295187289783
$ws.Range("J14").Value = 0.02712295187289783
$ws.Range("M14").Value = 4.993165333333334
$ws.Range("N14").Value = 14.979496
$ws.Range("O14").Value = 0.06779298131037136
$ws.Range("P14").Value = 0.06779298131037137
$ws.Range("Q14").Value = 0.3130315210773333
$ws.Range("R14").Value = 2.817283689696
$ws.Range("S14").Value = 0.001838745769401464
$ws.Range("T14").Value = 0.001838745769401465

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.062692
$ws.Range("H15").Value = 0.188076
$ws.Range("I15").Value = 0.02712295187289783
$ws.Range("J15").Value = 0.02712295187289783
$ws.Range("O15").Value = 0.5355771637189464
$ws.Range("P15").Value = 0.5355771637189464
$ws.Range("Q15").Value = 2.47300724902
$ws.Range("R15").Value = 22.25706524118
$ws.Range("S15").Value = 0.0145264336357721
$ws.Range("T15").Value = 0.0145264336357721

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.062692
$ws.Range("H16").Value = 0.188076
$ws.Range("I16").Value = 0.02712295187289783
$ws.Range("J16").Value = 0.02712295187289783
$ws.Range("M16").Value = 29.08216166666666
$ws.Range("N16").Value = 87.24648499999999
$ws.Range("O16").Value = 0.3948530262300277
$ws.Range("P16").Value = 0.3948530262300277
$ws.Range("Q16").Value = 1.823218879206666
$ws.Range("R16").Value = 16.40896991286
$ws.Range("S16").Value = 0.0107095796273051
$ws.Range("T16").Value = 0.01070957962730511

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.062692
$ws.Range("H17").Value = 0.188076
$ws.Range("I17").Value = 0.02712295187289783
$ws.Range("J17").Value = 0.02712295187289783
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.130869
$ws.Range("N17").Value = 0.392607
$ws.Range("O17").Value = 0.001776828740654623
$ws.Range("P17").Value = 0.001776828740654624
$ws.Range("Q17").Value = 0.008204439348000001
$ws.Range("R17").Value = 0.07383995413200001
$ws.Range("S17").Value = [double]"4.819284041915701E-05"
$ws.Range("T17").Value = [double]"4.819284041915702E-05"
